$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 5568.3335
$ws.Range("J121").Value = 4302.5
$ws.Range("L121").Value = 12907.5
$ws.Range("N121").Value = -16401.5

$ws.Range("H127").Value = 43487.875
$ws.Range("I127").Value = 2081.4285
$ws.Range("K127").Value = 6244.2855
$ws.Range("M127").Value = -1284.2855

$ws.Range("H132").Value = 5840.1226
$ws.Range("I132").Value = 3451.1843
$ws.Range("K132").Value = 10353.5529
$ws.Range("M132").Value = -7823.552899999999

$ws.Range("H137").Value = 5771.75
$ws.Range("I137").Value = 6609.625
$ws.Range("K137").Value = 19828.875
$ws.Range("M137").Value = -17278.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 986315.0600000001
$ws.Range("I2").Value = 1059168.1
$ws.Range("K2").Value = 1059168.1
$ws.Range("M2").Value = -1059055.1

$ws.Range("H63").Value = 3569.5625
$ws.Range("I63").Value = 3899
$ws.Range("J63").Value = 3558.9355
$ws.Range("K63").Value = 3899
$ws.Range("L63").Value = 3558.9355
$ws.Range("M63").Value = -3213
$ws.Range("N63").Value = -4930.9355

$ws.Range("H66").Value = 3569.5625
$ws.Range("I66").Value = 3899
$ws.Range("J66").Value = 3558.9355
$ws.Range("K66").Value = 19495
$ws.Range("L66").Value = 17794.6775
$ws.Range("M66").Value = -16063
$ws.Range("N66").Value = -24658.6775

$ws.Range("H116").Value = 986315.0600000001
$ws.Range("I116").Value = 1059168.1
$ws.Range("K116").Value = 1059168.1
$ws.Range("M116").Value = -1056874.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 986315.0600000001
$ws.Range("I3").Value = 1059168.1
$ws.Range("K3").Value = 1059168.1
$ws.Range("M3").Value = -1059054.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3294
$ws.Range("J58").Value = 3867.5
$ws.Range("L58").Value = 3867.5
$ws.Range("N58").Value = -4273.5

$ws.Range("H99").Value = 11695.634
$ws.Range("I99").Value = 7443.278
$ws.Range("K99").Value = 7443.278
$ws.Range("M99").Value = -5945.278

$ws.Range("H126").Value = 11695.634
$ws.Range("I126").Value = 7443.278
$ws.Range("K126").Value = 22329.834
$ws.Range("M126").Value = -19859.834

$ws.Range("H132").Value = 2376.1143
$ws.Range("I132").Value = 2270.8438
$ws.Range("J132").Value = 3499
$ws.Range("K132").Value = 6812.5314
$ws.Range("L132").Value = 10497
$ws.Range("M132").Value = -4282.5314
$ws.Range("N132").Value = -15557

$ws.Range("H136").Value = 3294
$ws.Range("J136").Value = 3867.5
$ws.Range("L136").Value = 11602.5
$ws.Range("N136").Value = -16702.5

$ws.Range("H141").Value = 129591.38
$ws.Range("J141").Value = 138113.1
$ws.Range("L141").Value = 138113.1
$ws.Range("N141").Value = -148473.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 954.06665
$ws.Range("I5").Value = 681.5
$ws.Range("K5").Value = 2044.5
$ws.Range("M5").Value = -1932.5

$ws.Range("H34").Value = 467.3
$ws.Range("I34").Value = 201
$ws.Range("J34").Value = 644.8333
$ws.Range("K34").Value = 603
$ws.Range("L34").Value = 1934.4999
$ws.Range("M34").Value = -519
$ws.Range("N34").Value = -2102.4999

$ws.Range("H39").Value = 8990.888999999999
$ws.Range("I39").Value = 800
$ws.Range("J39").Value = 9472.706
$ws.Range("K39").Value = 2400
$ws.Range("L39").Value = 28418.118
$ws.Range("M39").Value = -2106
$ws.Range("N39").Value = -29006.118

$ws.Range("H55").Value = 3085.8948
$ws.Range("J55").Value = 6525
$ws.Range("L55").Value = 19575
$ws.Range("N55").Value = -19929

$ws.Range("H70").Value = 4499
$ws.Range("J70").Value = 4999
$ws.Range("L70").Value = 14997
$ws.Range("N70").Value = -15627

$ws.Range("H73").Value = 4499
$ws.Range("J73").Value = 4999
$ws.Range("L73").Value = 14997
$ws.Range("N73").Value = -17181

$ws.Range("H74").Value = 3000
$ws.Range("I74").Value = 3000
$ws.Range("K74").Value = 9000
$ws.Range("M74").Value = -7939

$ws.Range("H77").Value = 3000
$ws.Range("I77").Value = 3000
$ws.Range("K77").Value = 27000
$ws.Range("M77").Value = -21696

$ws.Range("H80").Value = 2792.75
$ws.Range("I80").Value = 3002
$ws.Range("J80").Value = 2723
$ws.Range("K80").Value = 9006
$ws.Range("L80").Value = 8169
$ws.Range("M80").Value = -8070
$ws.Range("N80").Value = -10041

$ws.Range("H83").Value = 2792.75
$ws.Range("I83").Value = 3002
$ws.Range("J83").Value = 2723
$ws.Range("K83").Value = 27018
$ws.Range("L83").Value = 24507
$ws.Range("M83").Value = -22338
$ws.Range("N83").Value = -33867

$ws.Range("H120").Value = 885
$ws.Range("I120").Value = 885
$ws.Range("K120").Value = 2655
$ws.Range("M120").Value = 2183

$ws.Range("H134").Value = 3673.4666
$ws.Range("I134").Value = 3091.8333
$ws.Range("K134").Value = 9275.499899999999
$ws.Range("M134").Value = -4205.499899999999

$ws.Range("H135").Value = 954.06665
$ws.Range("I135").Value = 681.5
$ws.Range("K135").Value = 6133.5
$ws.Range("M135").Value = -3598.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4800.355
$ws.Range("I102").Value = 4941.9546
$ws.Range("J102").Value = 4454.222
$ws.Range("K102").Value = 4941.9546
$ws.Range("L102").Value = 4454.222
$ws.Range("M102").Value = -3319.9546
$ws.Range("N102").Value = -7698.222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3328.25
$ws.Range("I7").Value = 3026
$ws.Range("K7").Value = 3026
$ws.Range("M7").Value = -2914

$ws.Range("H46").Value = 2270.625
$ws.Range("I46").Value = 1045.2727
$ws.Range("J46").Value = 3307.4614
$ws.Range("K46").Value = 1045.2727
$ws.Range("L46").Value = 3307.4614
$ws.Range("M46").Value = -857.2727
$ws.Range("N46").Value = -3683.4614

$ws.Range("H126").Value = 3328.25
$ws.Range("I126").Value = 3026
$ws.Range("K126").Value = 9078
$ws.Range("M126").Value = -6608

$ws.Range("H132").Value = 55372.094
$ws.Range("I132").Value = 70001.07000000001
$ws.Range("J132").Value = 18799.666
$ws.Range("K132").Value = 210003.21
$ws.Range("L132").Value = 56398.99800000001
$ws.Range("M132").Value = -207473.21
$ws.Range("N132").Value = -61458.99800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 11502
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 11502
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 11502
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -14248
